$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column D rows 2-54 with the value "90minuteinduction"
$ws.Range("D2:D54").Value = "90minuteinduction"

# Update the saved selection to match the edited range
$ws.Range("D3:D54").Select()
